$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(119, 8).Value = 3000  # ALC!H119: 12500 -> 3000
$ws.Cells.Item(119, 10).Value = 1000  # ALC!J119: 15000 -> 1000
$ws.Cells.Item(119, 12).Value = 3000  # ALC!L119: 45000 -> 3000
$ws.Cells.Item(119, 14).Value = -12676  # ALC!N119: -54676 -> -12676

$ws.Cells.Item(129, 8).Value = 1071.75  # ALC!H129: 1016.6774 -> 1071.75
$ws.Cells.Item(129, 9).Value = 248.5  # ALC!I129: 265.66666 -> 248.5
$ws.Cells.Item(129, 10).Value = 1146.591  # ALC!J129: 1097.1428 -> 1146.591
$ws.Cells.Item(129, 11).Value = 745.5  # ALC!K129: 796.9999799999999 -> 745.5
$ws.Cells.Item(129, 12).Value = 3439.773  # ALC!L129: 3291.4284 -> 3439.773
$ws.Cells.Item(129, 13).Value = 4254.5  # ALC!M129: 4203.00002 -> 4254.5
$ws.Cells.Item(129, 14).Value = -13439.773  # ALC!N129: -13291.4284 -> -13439.773

$ws.Cells.Item(131, 8).Value = 10379.608  # ALC!H131: 11958.263 -> 10379.608
$ws.Cells.Item(131, 9).Value = 2980.4285  # ALC!I131: 3200.6365 -> 2980.4285
$ws.Cells.Item(131, 10).Value = 21889.445  # ALC!J131: 24000 -> 21889.445
$ws.Cells.Item(131, 11).Value = 8941.2855  # ALC!K131: 9601.9095 -> 8941.2855
$ws.Cells.Item(131, 12).Value = 65668.33499999999  # ALC!L131: 72000 -> 65668.33499999999
$ws.Cells.Item(131, 13).Value = -3901.2855  # ALC!M131: -4561.9095 -> -3901.2855
$ws.Cells.Item(131, 14).Value = -75748.33499999999  # ALC!N131: -82080 -> -75748.33499999999

$ws.Cells.Item(133, 8).Value = 49471.125  # ALC!H133: 49477.5 -> 49471.125
$ws.Cells.Item(133, 10).Value = 49471.125  # ALC!J133: 49477.5 -> 49471.125
$ws.Cells.Item(133, 12).Value = 49471.125  # ALC!L133: 49477.5 -> 49471.125
$ws.Cells.Item(133, 14).Value = -59591.125  # ALC!N133: -59597.5 -> -59591.125

$ws.Cells.Item(136, 8).Value = 63500  # ALC!H136: 0 -> 63500
$ws.Cells.Item(136, 10).Value = 63500  # ALC!J136: 0 -> 63500
$ws.Cells.Item(136, 12).Value = 63500  # ALC!L136: 0 -> 63500
$ws.Cells.Item(136, 14).Value = -73700  # ALC!N136: None -> -73700

$ws.Cells.Item(138, 8).Value = 6726210  # ALC!H138: 6283063.5 -> 6726210
$ws.Cells.Item(138, 9).Value = 2749653.8  # ALC!I138: 3249538.8 -> 2749653.8
$ws.Cells.Item(138, 10).Value = 8478591  # ALC!J138: 7250274.5 -> 8478591
$ws.Cells.Item(138, 11).Value = 8248961.399999999  # ALC!K138: 9748616.399999999 -> 8248961.399999999
$ws.Cells.Item(138, 12).Value = 25435773  # ALC!L138: 21750823.5 -> 25435773
$ws.Cells.Item(138, 13).Value = -8243821.399999999  # ALC!M138: -9743476.399999999 -> -8243821.399999999
$ws.Cells.Item(138, 14).Value = -25446053  # ALC!N138: -21761103.5 -> -25446053

$ws.Cells.Item(141, 8).Value = 2928.6667  # ALC!H141: 3194.5 -> 2928.6667
$ws.Cells.Item(141, 9).Value = 3014.4  # ALC!I141: 3216.111 -> 3014.4
$ws.Cells.Item(141, 10).Value = 2500  # ALC!J141: 3000 -> 2500
$ws.Cells.Item(141, 11).Value = 9043.200000000001  # ALC!K141: 9648.332999999999 -> 9043.200000000001
$ws.Cells.Item(141, 12).Value = 7500  # ALC!L141: 9000 -> 7500
$ws.Cells.Item(141, 13).Value = -3863.200000000001  # ALC!M141: -4468.332999999999 -> -3863.200000000001
$ws.Cells.Item(141, 14).Value = -17860  # ALC!N141: -19360 -> -17860

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 35539.7  # ARM!H32: 36097.387 -> 35539.7
$ws.Cells.Item(32, 9).Value = 9358.154  # ARM!I32: 9388.289000000001 -> 9358.154
$ws.Cells.Item(32, 11).Value = 9358.154  # ARM!K32: 9388.289000000001 -> 9358.154
$ws.Cells.Item(32, 13).Value = -9071.154  # ARM!M32: -9101.289000000001 -> -9071.154

$ws.Cells.Item(132, 8).Value = 2148.38  # ARM!H132: 1947.5518 -> 2148.38
$ws.Cells.Item(132, 9).Value = 1888.7333  # ARM!I132: 1819.7551 -> 1888.7333
$ws.Cells.Item(132, 10).Value = 4485.2  # ARM!J132: 2643.3333 -> 4485.2
$ws.Cells.Item(132, 11).Value = 5666.199900000001  # ARM!K132: 5459.2653 -> 5666.199900000001
$ws.Cells.Item(132, 12).Value = 13455.6  # ARM!L132: 7929.999899999999 -> 13455.6
$ws.Cells.Item(132, 13).Value = -3136.199900000001  # ARM!M132: -2929.2653 -> -3136.199900000001
$ws.Cells.Item(132, 14).Value = -18515.6  # ARM!N132: -12989.9999 -> -18515.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(122, 8).Value = 0  # BSM!H122: 33000 -> 0
$ws.Cells.Item(122, 10).Value = 0  # BSM!J122: 33000 -> 0
$ws.Cells.Item(122, 12).Value = 0  # BSM!L122: 33000 -> 0
$ws.Cells.Item(122, 14).ClearContents()  # BSM!N122: -42800 -> (removed)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1612.9445  # CRP!H16: 1891.4 -> 1612.9445
$ws.Cells.Item(16, 9).Value = 1472.3077  # CRP!I16: 1804.2 -> 1472.3077
$ws.Cells.Item(16, 11).Value = 1472.3077  # CRP!K16: 1804.2 -> 1472.3077
$ws.Cells.Item(16, 13).Value = -1185.3077  # CRP!M16: -1517.2 -> -1185.3077

$ws.Cells.Item(64, 8).Value = 20000  # CRP!H64: 0 -> 20000
$ws.Cells.Item(64, 10).Value = 20000  # CRP!J64: 0 -> 20000
$ws.Cells.Item(64, 12).Value = 20000  # CRP!L64: 0 -> 20000
$ws.Cells.Item(64, 14).Value = -20496  # CRP!N64: None -> -20496

$ws.Cells.Item(67, 8).Value = 20000  # CRP!H67: 0 -> 20000
$ws.Cells.Item(67, 10).Value = 20000  # CRP!J67: 0 -> 20000
$ws.Cells.Item(67, 12).Value = 20000  # CRP!L67: 0 -> 20000
$ws.Cells.Item(67, 14).Value = -21716  # CRP!N67: None -> -21716

$ws.Cells.Item(93, 8).Value = 14256.571  # CRP!H93: 13447.875 -> 14256.571
$ws.Cells.Item(93, 9).Value = 14256.571  # CRP!I93: 13447.875 -> 14256.571
$ws.Cells.Item(93, 11).Value = 14256.571  # CRP!K93: 13447.875 -> 14256.571
$ws.Cells.Item(93, 13).Value = -12384.571  # CRP!M93: -11575.875 -> -12384.571

$ws.Cells.Item(113, 8).Value = 1612.9445  # CRP!H113: 1891.4 -> 1612.9445
$ws.Cells.Item(113, 9).Value = 1472.3077  # CRP!I113: 1804.2 -> 1472.3077
$ws.Cells.Item(113, 11).Value = 1472.3077  # CRP!K113: 1804.2 -> 1472.3077
$ws.Cells.Item(113, 13).Value = 697.6922999999999  # CRP!M113: 365.8 -> 697.6922999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1940.3334  # CUL!H5: 2021.5807 -> 1940.3334
$ws.Cells.Item(5, 9).Value = 899.3889  # CUL!I5: 921.7646999999999 -> 899.3889
$ws.Cells.Item(5, 10).Value = 3189.4666  # CUL!J5: 3357.0715 -> 3189.4666
$ws.Cells.Item(5, 11).Value = 2698.1667  # CUL!K5: 2765.2941 -> 2698.1667
$ws.Cells.Item(5, 12).Value = 9568.399800000001  # CUL!L5: 10071.2145 -> 9568.399800000001
$ws.Cells.Item(5, 13).Value = -2586.1667  # CUL!M5: -2653.2941 -> -2586.1667
$ws.Cells.Item(5, 14).Value = -9792.399800000001  # CUL!N5: -10295.2145 -> -9792.399800000001

$ws.Cells.Item(8, 8).Value = 866.75  # CUL!H8: 40.88889 -> 866.75
$ws.Cells.Item(8, 9).Value = 866.75  # CUL!I8: 40.88889 -> 866.75
$ws.Cells.Item(8, 11).Value = 2600.25  # CUL!K8: 122.66667 -> 2600.25
$ws.Cells.Item(8, 13).Value = -2461.25  # CUL!M8: 16.33332999999999 -> -2461.25

$ws.Cells.Item(34, 8).Value = 2200  # CUL!H34: 2800 -> 2200

$ws.Cells.Item(39, 8).Value = 9250.15  # CUL!H39: 9263.315000000001 -> 9250.15
$ws.Cells.Item(39, 10).Value = 9684.421  # CUL!J39: 9722.444 -> 9684.421
$ws.Cells.Item(39, 12).Value = 29053.263  # CUL!L39: 29167.332 -> 29053.263
$ws.Cells.Item(39, 14).Value = -29641.263  # CUL!N39: -29755.332 -> -29641.263

$ws.Cells.Item(55, 8).Value = 3600  # CUL!H55: 2600 -> 3600
$ws.Cells.Item(55, 10).Value = 4857.143  # CUL!J55: 5500 -> 4857.143
$ws.Cells.Item(55, 12).Value = 14571.429  # CUL!L55: 16500 -> 14571.429
$ws.Cells.Item(55, 14).Value = -14925.429  # CUL!N55: -16854 -> -14925.429

$ws.Cells.Item(131, 8).Value = 6804335.5  # CUL!H131: 5748651.5 -> 6804335.5
$ws.Cells.Item(131, 10).Value = 6946074  # CUL!J131: 5849489 -> 6946074
$ws.Cells.Item(131, 12).Value = 20838222  # CUL!L131: 17548467 -> 20838222
$ws.Cells.Item(131, 14).Value = -20848302  # CUL!N131: -17558547 -> -20848302

$ws.Cells.Item(132, 8).Value = 1445  # CUL!H132: 1479.5 -> 1445
$ws.Cells.Item(132, 9).Value = 900  # CUL!I132: 1278 -> 900
$ws.Cells.Item(132, 10).Value = 1505.5555  # CUL!J132: 1546.6666 -> 1505.5555
$ws.Cells.Item(132, 11).Value = 8100  # CUL!K132: 11502 -> 8100
$ws.Cells.Item(132, 12).Value = 13549.9995  # CUL!L132: 13919.9994 -> 13549.9995
$ws.Cells.Item(132, 13).Value = -5570  # CUL!M132: -8972 -> -5570
$ws.Cells.Item(132, 14).Value = -18609.9995  # CUL!N132: -18979.9994 -> -18609.9995

$ws.Cells.Item(135, 8).Value = 1940.3334  # CUL!H135: 2021.5807 -> 1940.3334
$ws.Cells.Item(135, 9).Value = 899.3889  # CUL!I135: 921.7646999999999 -> 899.3889
$ws.Cells.Item(135, 10).Value = 3189.4666  # CUL!J135: 3357.0715 -> 3189.4666
$ws.Cells.Item(135, 11).Value = 8094.5001  # CUL!K135: 8295.882299999999 -> 8094.5001
$ws.Cells.Item(135, 12).Value = 28705.1994  # CUL!L135: 30213.6435 -> 28705.1994
$ws.Cells.Item(135, 13).Value = -5559.5001  # CUL!M135: -5760.882299999999 -> -5559.5001
$ws.Cells.Item(135, 14).Value = -33775.1994  # CUL!N135: -35283.6435 -> -33775.1994

$ws.Cells.Item(137, 8).Value = 10104958  # CUL!H137: 4952 -> 10104958
$ws.Cells.Item(137, 9).Value = 16672009  # CUL!I137: 5377.5 -> 16672009
$ws.Cells.Item(137, 10).Value = 254383.25  # CUL!J137: 3250 -> 254383.25
$ws.Cells.Item(137, 11).Value = 50016027  # CUL!K137: 16132.5 -> 50016027
$ws.Cells.Item(137, 12).Value = 763149.75  # CUL!L137: 9750 -> 763149.75
$ws.Cells.Item(137, 13).Value = -50010927  # CUL!M137: -11032.5 -> -50010927
$ws.Cells.Item(137, 14).Value = -773349.75  # CUL!N137: -19950 -> -773349.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 45457230  # GSM!H80: 65219990 -> 45457230
$ws.Cells.Item(80, 9).Value = 2680.3572  # GSM!I80: 2675 -> 2680.3572
$ws.Cells.Item(80, 10).Value = 125002700  # GSM!J80: 187502460 -> 125002700
$ws.Cells.Item(80, 11).Value = 2680.3572  # GSM!K80: 2675 -> 2680.3572
$ws.Cells.Item(80, 12).Value = 125002700  # GSM!L80: 187502460 -> 125002700
$ws.Cells.Item(80, 13).Value = -1682.3572  # GSM!M80: -1677 -> -1682.3572
$ws.Cells.Item(80, 14).Value = -125004696  # GSM!N80: -187504456 -> -125004696

$ws.Cells.Item(83, 8).Value = 45457230  # GSM!H83: 65219990 -> 45457230
$ws.Cells.Item(83, 9).Value = 2680.3572  # GSM!I83: 2675 -> 2680.3572
$ws.Cells.Item(83, 10).Value = 125002700  # GSM!J83: 187502460 -> 125002700
$ws.Cells.Item(83, 11).Value = 13401.786  # GSM!K83: 13375 -> 13401.786
$ws.Cells.Item(83, 12).Value = 625013500  # GSM!L83: 937512300 -> 625013500
$ws.Cells.Item(83, 13).Value = -8409.786  # GSM!M83: -8383 -> -8409.786
$ws.Cells.Item(83, 14).Value = -625023484  # GSM!N83: -937522284 -> -625023484

$ws.Cells.Item(126, 8).Value = 2648.617  # GSM!H126: 2642.1914 -> 2648.617
$ws.Cells.Item(126, 9).Value = 2366.8948  # GSM!I126: 2351 -> 2366.8948
$ws.Cells.Item(126, 11).Value = 7100.6844  # GSM!K126: 7053 -> 7100.6844
$ws.Cells.Item(126, 13).Value = -4630.6844  # GSM!M126: -4583 -> -4630.6844

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 10).Value = 3000  # LTW!J20: 0 -> 3000
$ws.Cells.Item(20, 12).Value = 3000  # LTW!L20: 0 -> 3000
$ws.Cells.Item(20, 14).Value = -3452  # LTW!N20: None -> -3452

$ws.Cells.Item(22, 8).Value = 7075.625  # LTW!H22: 8050.7144 -> 7075.625
$ws.Cells.Item(22, 10).Value = 14697.143  # LTW!J22: 20476 -> 14697.143
$ws.Cells.Item(22, 12).Value = 14697.143  # LTW!L22: 20476 -> 14697.143
$ws.Cells.Item(22, 14).Value = -15287.143  # LTW!N22: -21066 -> -15287.143

$ws.Cells.Item(27, 8).Value = 7075.625  # LTW!H27: 8050.7144 -> 7075.625
$ws.Cells.Item(27, 10).Value = 14697.143  # LTW!J27: 20476 -> 14697.143
$ws.Cells.Item(27, 12).Value = 14697.143  # LTW!L27: 20476 -> 14697.143
$ws.Cells.Item(27, 14).Value = -14911.143  # LTW!N27: -20690 -> -14911.143

$ws.Cells.Item(64, 8).Value = 18996.666  # LTW!H64: 19396 -> 18996.666
$ws.Cells.Item(64, 10).Value = 18996.666  # LTW!J64: 19396 -> 18996.666
$ws.Cells.Item(64, 12).Value = 18996.666  # LTW!L64: 19396 -> 18996.666
$ws.Cells.Item(64, 14).Value = -19446.666  # LTW!N64: -19846 -> -19446.666

$ws.Cells.Item(67, 8).Value = 18996.666  # LTW!H67: 19396 -> 18996.666
$ws.Cells.Item(67, 10).Value = 18996.666  # LTW!J67: 19396 -> 18996.666
$ws.Cells.Item(67, 12).Value = 18996.666  # LTW!L67: 19396 -> 18996.666
$ws.Cells.Item(67, 14).Value = -20556.666  # LTW!N67: -20956 -> -20556.666

$ws.Cells.Item(69, 8).Value = 30000  # LTW!H69: 0 -> 30000
$ws.Cells.Item(69, 10).Value = 30000  # LTW!J69: 0 -> 30000
$ws.Cells.Item(69, 12).Value = 30000  # LTW!L69: 0 -> 30000
$ws.Cells.Item(69, 14).Value = -31622  # LTW!N69: None -> -31622

$ws.Cells.Item(72, 8).Value = 30000  # LTW!H72: 0 -> 30000
$ws.Cells.Item(72, 10).Value = 30000  # LTW!J72: 0 -> 30000
$ws.Cells.Item(72, 12).Value = 90000  # LTW!L72: 0 -> 90000
$ws.Cells.Item(72, 14).Value = -98112  # LTW!N72: None -> -98112

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 37599.8  # WVR!H62: 57000 -> 37599.8
$ws.Cells.Item(62, 9).Value = 57500  # WVR!I62: 95000 -> 57500
$ws.Cells.Item(62, 10).Value = 24333  # WVR!J62: 19000 -> 24333
$ws.Cells.Item(62, 11).Value = 57500  # WVR!K62: 95000 -> 57500
$ws.Cells.Item(62, 12).Value = 24333  # WVR!L62: 19000 -> 24333
$ws.Cells.Item(62, 13).Value = -56876  # WVR!M62: -94376 -> -56876
$ws.Cells.Item(62, 14).Value = -25581  # WVR!N62: -20248 -> -25581

$ws.Cells.Item(64, 8).Value = 29038  # WVR!H64: 28052.334 -> 29038
$ws.Cells.Item(64, 9).Value = 0  # WVR!I64: 28000 -> 0
$ws.Cells.Item(64, 10).Value = 29038  # WVR!J64: 28062.8 -> 29038
$ws.Cells.Item(64, 11).Value = 0  # WVR!K64: 28000 -> 0
$ws.Cells.Item(64, 12).Value = 29038  # WVR!L64: 28062.8 -> 29038
$ws.Cells.Item(64, 13).ClearContents()  # WVR!M64: -27752 -> (removed)
$ws.Cells.Item(64, 14).Value = -29534  # WVR!N64: -28558.8 -> -29534

$ws.Cells.Item(65, 8).Value = 37599.8  # WVR!H65: 57000 -> 37599.8
$ws.Cells.Item(65, 9).Value = 57500  # WVR!I65: 95000 -> 57500
$ws.Cells.Item(65, 10).Value = 24333  # WVR!J65: 19000 -> 24333
$ws.Cells.Item(65, 11).Value = 287500  # WVR!K65: 475000 -> 287500
$ws.Cells.Item(65, 12).Value = 121665  # WVR!L65: 95000 -> 121665
$ws.Cells.Item(65, 13).Value = -284380  # WVR!M65: -471880 -> -284380
$ws.Cells.Item(65, 14).Value = -127905  # WVR!N65: -101240 -> -127905

$ws.Cells.Item(67, 8).Value = 29038  # WVR!H67: 28052.334 -> 29038
$ws.Cells.Item(67, 9).Value = 0  # WVR!I67: 28000 -> 0
$ws.Cells.Item(67, 10).Value = 29038  # WVR!J67: 28062.8 -> 29038
$ws.Cells.Item(67, 11).Value = 0  # WVR!K67: 28000 -> 0
$ws.Cells.Item(67, 12).Value = 29038  # WVR!L67: 28062.8 -> 29038
$ws.Cells.Item(67, 13).ClearContents()  # WVR!M67: -27142 -> (removed)
$ws.Cells.Item(67, 14).Value = -30754  # WVR!N67: -29778.8 -> -30754

$ws.Cells.Item(81, 8).Value = 4137.1665  # WVR!H81: 1024.3636 -> 4137.1665
$ws.Cells.Item(81, 9).Value = 860  # WVR!I81: 726.8 -> 860
$ws.Cells.Item(81, 10).Value = 4999.579  # WVR!J81: 4000 -> 4999.579
$ws.Cells.Item(81, 11).Value = 1720  # WVR!K81: 1453.6 -> 1720
$ws.Cells.Item(81, 12).Value = 9999.157999999999  # WVR!L81: 8000 -> 9999.157999999999
$ws.Cells.Item(81, 13).Value = -659  # WVR!M81: -392.5999999999999 -> -659
$ws.Cells.Item(81, 14).Value = -12121.158  # WVR!N81: -10122 -> -12121.158

$ws.Cells.Item(84, 8).Value = 4137.1665  # WVR!H84: 1024.3636 -> 4137.1665
$ws.Cells.Item(84, 9).Value = 860  # WVR!I84: 726.8 -> 860
$ws.Cells.Item(84, 10).Value = 4999.579  # WVR!J84: 4000 -> 4999.579
$ws.Cells.Item(84, 11).Value = 8600  # WVR!K84: 7268 -> 8600
$ws.Cells.Item(84, 12).Value = 49995.78999999999  # WVR!L84: 40000 -> 49995.78999999999
$ws.Cells.Item(84, 13).Value = -3296  # WVR!M84: -1964 -> -3296
$ws.Cells.Item(84, 14).Value = -60603.78999999999  # WVR!N84: -50608 -> -60603.78999999999

$ws.Cells.Item(100, 8).Value = 688.5714  # WVR!H100: 588.7273 -> 688.5714
$ws.Cells.Item(100, 9).Value = 736.6667  # WVR!I100: 588.7273 -> 736.6667
$ws.Cells.Item(100, 10).Value = 400  # WVR!J100: 0 -> 400
$ws.Cells.Item(100, 11).Value = 1473.3334  # WVR!K100: 1177.4546 -> 1473.3334
$ws.Cells.Item(100, 12).Value = 800  # WVR!L100: 0 -> 800
$ws.Cells.Item(100, 13).Value = -932.3334  # WVR!M100: -636.4546 -> -932.3334
$ws.Cells.Item(100, 14).Value = -1882  # WVR!N100: None -> -1882

$ws.Cells.Item(122, 8).Value = 1266.6666  # WVR!H122: 1325 -> 1266.6666
$ws.Cells.Item(122, 9).Value = 1300  # WVR!I122: 1325 -> 1300
$ws.Cells.Item(122, 10).Value = 1200  # WVR!J122: 0 -> 1200
$ws.Cells.Item(122, 11).Value = 3900  # WVR!K122: 3975 -> 3900
$ws.Cells.Item(122, 12).Value = 3600  # WVR!L122: 0 -> 3600
$ws.Cells.Item(122, 13).Value = -1450  # WVR!M122: -1525 -> -1450
$ws.Cells.Item(122, 14).Value = -8500  # WVR!N122: None -> -8500
